# Update the "想去人数" (want-to-go count) values for two exhibition
# entries that appear in both the "展览" sheet and the aggregated
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet: rows 3 & 4, column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1216
$wsExhibition.Range("F4").Value = 2691

# 全部类型 (All Types) sheet: rows 5 & 6, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1216
$wsAll.Range("F6").Value = 2691
